# Wrote a new test case: bump the "Or Higher" test-case count for the
# AppendOverlay row (row 3) from 4 to 5. Excel recalculates the dependent
# summary formulas (G4 = SUM($C:$C), G6 = G5/G4) automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = 5

# Leave the selection where the next test case would be entered.
$ws.Range("E3").Select()
